$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save the current values of rows 9, 14, 15 (A:C) before overwriting
$row9  = @($ws.Range("A9").Value2,  $ws.Range("B9").Value2,  $ws.Range("C9").Value2)
$row14 = @($ws.Range("A14").Value2, $ws.Range("B14").Value2, $ws.Range("C14").Value2)
$row15 = @($ws.Range("A15").Value2, $ws.Range("B15").Value2, $ws.Range("C15").Value2)

# Rotate: row9 <- old row15, row14 <- old row9, row15 <- old row14
$ws.Range("A9").Value2  = $row15[0]
$ws.Range("B9").Value2  = $row15[1]
$ws.Range("C9").Value2  = $row15[2]

$ws.Range("A14").Value2 = $row9[0]
$ws.Range("B14").Value2 = $row9[1]
$ws.Range("C14").Value2 = $row9[2]

$ws.Range("A15").Value2 = $row14[0]
$ws.Range("B15").Value2 = $row14[1]
$ws.Range("C15").Value2 = $row14[2]
